$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "No Pic, sorry :("

$ws.Range("C14").Select()
